$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look like plain numbers need the Text number format
# forced first, so Excel keeps storing them as text (matching the source data,
# which uses text-formatted price columns, e.g. "2.194.40").
$ws.Range("D2").Value = '35.547.85'
$ws.Range("E2").Value = '  +1.46%  '
$ws.Range("D3").Value = '1.911.61'
$ws.Range("E3").Value = '  +3.15%  '
$ws.Range("E4").Value = '  +0.50%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '246.16'
$ws.Range("E5").Value = '  +3.85%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.655'
$ws.Range("E6").Value = '  +5.27%  '
$ws.Range("E7").Value = '  +0.50%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '41.96'
$ws.Range("E8").Value = '  -0.99%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.346'
$ws.Range("E9").Value = '  +5.47%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '49.45'
$ws.Range("E10").Value = '  +5.71%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0719'
$ws.Range("E11").Value = '  +3.54%  '
$ws.Range("E12").Value = '  +0.76%  '
$ws.Range("D13").Value = '2.190.02'
$ws.Range("E13").Value = '  +3.24%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '12.27'
$ws.Range("E14").Value = '  +7.32%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.699'
$ws.Range("E15").Value = '  +3.41%  '
$ws.Range("B16").Value = 'Polkadot'
$ws.Range("C16").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '4.89'
$ws.Range("E16").Value = '  +1.73%  '
$ws.Range("B17").Value = 'WrappedEther'
$ws.Range("C17").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D17").Value = '1.894.56'
$ws.Range("E17").Value = '  +2.08%  '
$ws.Range("D18").Value = '35.532.35'
$ws.Range("E18").Value = '  +1.41%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '72.41'
$ws.Range("E19").Value = '  +3.02%  '
$ws.Range("D20").Value = '0.0₃0822'
$ws.Range("E20").Value = '  +3.27%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '245.30'
$ws.Range("E21").Value = '  +1.99%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '12.63'
$ws.Range("E22").Value = '  +3.60%  '
$ws.Range("E23").Value = '  +1.24%  '
$ws.Range("E24").Value = '  +0.48%  '
$ws.Range("E25").Value = '  +1.67%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.21'
$ws.Range("E26").Value = '  +17.97%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '171.38'
$ws.Range("E27").Value = '  +0.48%  '
$ws.Range("E28").Value = '  +5.27%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '18.43'
$ws.Range("E29").Value = '  +4.36%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.128'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.16'
$ws.Range("E31").Value = '  +3.69%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.0570'
$ws.Range("E32").Value = '  +2.14%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.02'
$ws.Range("E33").Value = '  +0.55%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.15'
$ws.Range("E34").Value = '  +2.98%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.925'
$ws.Range("E35").Value = '  +18.54%  '
$ws.Range("E36").Value = '  +5.15%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.03'
$ws.Range("E37").Value = '  +1.46%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.34'
$ws.Range("E38").Value = '  +0.67%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0211'
$ws.Range("E39").Value = '  +4.39%  '
$ws.Range("E40").Value = '  +2.44%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0639'
$ws.Range("E41").Value = '  +15.52%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '91.21'
$ws.Range("E42").Value = '  +0.50%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '15.78'
$ws.Range("E43").Value = '  +7.28%  '
$ws.Range("D44").Value = '1.353.56'
$ws.Range("E44").Value = '  +0.17%  '
$ws.Range("E45").Value = '  +3.15%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '47.90'
$ws.Range("E46").Value = '  +38.51%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '12.70'
$ws.Range("E47").Value = '  -0.69%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.79'
$ws.Range("E48").Value = '  +1.97%  '
$ws.Range("E49").Value = '  -0.06%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '6.56'
$ws.Range("E50").Value = '  +0.32%  '
$ws.Range("D51").Value = '2.098.01'
$ws.Range("E51").Value = '  +3.16%  '
